$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.053.30"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "2.636.47"
$ws.Range("E3").Value = "  +3.84%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'518.29"
$ws.Range("E5").Value = "  +2.36%  "
$ws.Range("D6").Value = "'145.79"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").Value = "'0.568"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").Value = "2.661.90"
$ws.Range("E9").Value = "  +4.69%  "
$ws.Range("D10").Value = "'6.27"
$ws.Range("E10").Value = "  +2.95%  "
$ws.Range("E11").Value = "  +2.99%  "
$ws.Range("D12").Value = "'0.338"
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").Value = "3.104.57"
$ws.Range("E14").Value = "  +4.08%  "
$ws.Range("D15").Value = "58.999.70"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "'21.03"
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").Value = "'0.0000138"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("D18").Value = "2.656.84"
$ws.Range("E18").Value = "  +4.60%  "
$ws.Range("D19").Value = "'350.39"
$ws.Range("E19").Value = "  +3.39%  "
$ws.Range("D20").Value = "'4.53"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "'10.38"
$ws.Range("E21").Value = "  +2.82%  "
$ws.Range("E22").Value = "  +4.15%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'61.86"
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("D25").Value = "'0.421"
$ws.Range("E25").Value = "  +2.38%  "
$ws.Range("D26").Value = "2.765.57"
$ws.Range("E26").Value = "  +4.30%  "
$ws.Range("D27").Value = "'0.163"
$ws.Range("E27").Value = "  +2.00%  "
$ws.Range("D28").Value = "'0.995"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("D29").Value = "0.0₃0809"
$ws.Range("E29").Value = "  +2.83%  "
$ws.Range("D30").Value = "'7.16"
$ws.Range("E30").Value = "  +3.05%  "
$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").Value = "'6.28"
$ws.Range("E32").Value = "  +7.80%  "
$ws.Range("D33").Value = "'19.03"
$ws.Range("E33").Value = "  +2.66%  "
$ws.Range("E34").Value = "  +3.08%  "
$ws.Range("D35").Value = "'150.22"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").Value = "'0.966"
$ws.Range("E36").Value = "  +5.73%  "
$ws.Range("D37").Value = "'4.03"
$ws.Range("E37").Value = "  +3.51%  "
$ws.Range("E38").Value = "  +2.71%  "
$ws.Range("D39").Value = "'36.80"
$ws.Range("E39").Value = "  +1.96%  "
$ws.Range("D40").Value = "'0.847"
$ws.Range("E40").Value = "  +3.03%  "
$ws.Range("D41").Value = "'3.72"
$ws.Range("E41").Value = "  +5.50%  "
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("D43").Value = "'278.19"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.610"
$ws.Range("E44").Value = "  +1.70%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'0.994"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").Value = "'0.0985"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("D47").Value = "'19.63"
$ws.Range("E47").Value = "  +5.20%  "
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("D49").Value = "'10.30"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'4.72"
$ws.Range("E50").Value = "  +4.50%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0230"
$ws.Range("E51").Value = "  +1.57%  "
